# Refresh the "cryptos" price table (GitHub Actions scheduled update).
# Only the Coin/Link/Price/Volume(1h) cells that actually changed between
# runs are touched; row 1 (headers) and column A (rank index) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data keeps every "Price" cell as literal text (e.g. "37.250.90" or
# "0.388"), never a real number -- that's how the scraper writes the sheet,
# and it's also what lets a thousands-grouped value like "37.250.90" round
# -trip unchanged. Writing a plain numeric-looking string through
# Range.Value lets Excel's normal type inference silently coerce it to a
# number (and drop the trailing zero, e.g. "247.16" -> 247.16, "1.00" -> 1),
# so values that parse as a plain decimal get a leading apostrophe -- the
# same "force text" quote-prefix trick you'd use typing into a cell by hand.
function Set-TextValue {
    param($range, [string]$value)

    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($range).Value = "'" + $value
    } else {
        $ws.Range($range).Value = $value
    }
}

# Row 2 - Bitcoin
Set-TextValue "D2" "37.219.21"
$ws.Range("E2").Value = "  +1.39%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.016.30"
$ws.Range("E3").Value = "  +2.68%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "247.16"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.50%  "

# Row 7 - Solana
Set-TextValue "D7" "59.93"
$ws.Range("E7").Value = "  -2.27%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.09%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.388"
$ws.Range("E9").Value = "  +2.68%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0809"
$ws.Range("E10").Value = "  +1.55%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.18%  "

# Row 12 - Chainlink
Set-TextValue "D12" "15.17"
$ws.Range("E12").Value = "  +5.63%  "

# Row 13 - was Avalanche, now WrappedliquidstakedEther2.0 (rows 13-15 reorder)
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D13" "2.315.78"
$ws.Range("E13").Value = "  +3.67%  "

# Row 14 - was WrappedliquidstakedEther2.0, now Polygon
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D14" "0.854"
$ws.Range("E14").Value = "  +1.44%  "

# Row 15 - was Polygon, now Avalanche
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D15" "22.32"
$ws.Range("E15").Value = "  +0.81%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.47"
$ws.Range("E16").Value = "  +2.90%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.025.43"
$ws.Range("E17").Value = "  +2.89%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "37.143.63"
$ws.Range("E18").Value = "  +1.49%  "

# Row 19 - Litecoin
Set-TextValue "D19" "70.46"
$ws.Range("E19").Value = "  +0.63%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +1.16%  "

# Row 21 - Uniswap
Set-TextValue "D21" "5.23"
$ws.Range("E21").Value = "  +2.43%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "231.16"
$ws.Range("E22").Value = "  +0.17%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.04%  "

# Row 24 - PancakeSwap
Set-TextValue "D24" "2.49"
$ws.Range("E24").Value = "  +0.74%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +0.50%  "

# Row 26 - Cosmos
Set-TextValue "D26" "9.46"
$ws.Range("E26").Value = "  +2.24%  "

# Row 27 - Monero
Set-TextValue "D27" "164.28"
$ws.Range("E27").Value = "  +2.13%  "

# Row 28 - Kaspa
Set-TextValue "D28" "0.138"
$ws.Range("E28").Value = "  -4.29%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "19.79"
$ws.Range("E29").Value = "  +1.74%  "

# Row 30 - ImmutableX
Set-TextValue "D30" "1.39"
$ws.Range("E30").Value = "  +13.03%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.122"
$ws.Range("E31").Value = "  +1.57%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.32%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0663"
$ws.Range("E33").Value = "  +6.97%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -0.19%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "2.47"
$ws.Range("E35").Value = "  +8.00%  "

# Row 36 - was RenderToken, now BinanceUSD (rows 36-37 reorder)
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  +0.19%  "

# Row 37 - was BinanceUSD, now RenderToken
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D37" "3.43"
$ws.Range("E37").Value = "  -3.38%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  +2.20%  "

# Row 39 - THORChain
Set-TextValue "D39" "5.40"
$ws.Range("E39").Value = "  -3.62%  "

# Row 40 - Cronos
$ws.Range("E40").Value = "  -0.66%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +0.61%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "1.19"
$ws.Range("E42").Value = "  +0.96%  "

# Row 43 - VeChain
Set-TextValue "D43" "0.0215"
$ws.Range("E43").Value = "  +1.38%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "16.67"
$ws.Range("E44").Value = "  +2.31%  "

# Row 45 - Aave
Set-TextValue "D45" "92.27"
$ws.Range("E45").Value = "  +3.64%  "

# Row 46 - was ARBITRUM, now Maker (rows 46-47 reorder)
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D46" "1.381.25"
$ws.Range("E46").Value = "  +0.97%  "

# Row 47 - was Maker, now ARBITRUM
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D47" "1.06"
$ws.Range("E47").Value = "  +2.01%  "

# Row 48 - FraxShare
Set-TextValue "D48" "7.49"
$ws.Range("E48").Value = "  +4.51%  "

# Row 49 - NEARProtocol
$ws.Range("E49").Value = "  +12.72%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +0.15%  "

# Row 51 - MultiversX
Set-TextValue "D51" "46.80"
$ws.Range("E51").Value = "  +5.17%  "
